$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" column (C) date serial values for rows 2-5
# from 45184 to 45185, matching the target diff.
$ws.Range("C2:C5").Value = 45185
